$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale annotation cell that does not belong to the relocated record
$ws.Cells.Item(100, 6).Value = $null

# Row 2: A 43326-2025
$ws.Cells.Item(2, 1).Value = "A 43326-2025"
$ws.Cells.Item(2, 2).Value = 45910
$ws.Cells.Item(2, 3).Value = 46081
$ws.Cells.Item(2, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(2, 5).Value = "SVALÖV"
$ws.Cells.Item(2, 7).Value = 2.1
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 4
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 4
$ws.Cells.Item(2, 18).Value = "Entita`r`nGulsparv`r`nMindre hackspett`r`nSpillkråka"
$ws.Cells.Item(2, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 43326-2025 artfynd.xlsx", "A 43326-2025")'
$ws.Cells.Item(2, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 43326-2025 karta.png", "A 43326-2025")'
$ws.Cells.Item(2, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 43326-2025 FSC-klagomål.docx", "A 43326-2025")'
$ws.Cells.Item(2, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 43326-2025 FSC-klagomål mail.docx", "A 43326-2025")'
$ws.Cells.Item(2, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 43326-2025 tillsynsbegäran.docx", "A 43326-2025")'
$ws.Cells.Item(2, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 43326-2025 tillsynsbegäran mail.docx", "A 43326-2025")'
$ws.Cells.Item(2, 26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/fåglar/A 43326-2025 prioriterade fågelarter.docx", "A 43326-2025")'

# Row 3: A 39165-2024
$ws.Cells.Item(3, 1).Value = "A 39165-2024"
$ws.Cells.Item(3, 2).Value = 45548
$ws.Cells.Item(3, 3).Value = 46081
$ws.Cells.Item(3, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(3, 5).Value = "SVALÖV"
$ws.Cells.Item(3, 7).Value = 2.9
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 2
$ws.Cells.Item(3, 16).Value = 2
$ws.Cells.Item(3, 17).Value = 2
$ws.Cells.Item(3, 18).Value = "Skogsalm`r`nAsk"
$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 39165-2024 artfynd.xlsx", "A 39165-2024")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 39165-2024 karta.png", "A 39165-2024")'
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 39165-2024 FSC-klagomål.docx", "A 39165-2024")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 39165-2024 FSC-klagomål mail.docx", "A 39165-2024")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 39165-2024 tillsynsbegäran.docx", "A 39165-2024")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 39165-2024 tillsynsbegäran mail.docx", "A 39165-2024")'

# Row 4: A 38522-2025
$ws.Cells.Item(4, 1).Value = "A 38522-2025"
$ws.Cells.Item(4, 2).Value = 45884
$ws.Cells.Item(4, 3).Value = 46081
$ws.Cells.Item(4, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(4, 5).Value = "SVALÖV"
$ws.Cells.Item(4, 7).Value = 2.6
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = "Skogsknipprot`r`nStor häxört"
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 38522-2025 artfynd.xlsx", "A 38522-2025")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 38522-2025 karta.png", "A 38522-2025")'
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 38522-2025 FSC-klagomål.docx", "A 38522-2025")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 38522-2025 FSC-klagomål mail.docx", "A 38522-2025")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 38522-2025 tillsynsbegäran.docx", "A 38522-2025")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 38522-2025 tillsynsbegäran mail.docx", "A 38522-2025")'

# Row 5: A 61064-2024
$ws.Cells.Item(5, 1).Value = "A 61064-2024"
$ws.Cells.Item(5, 2).Value = 45645.49443287037
$ws.Cells.Item(5, 3).Value = 46081
$ws.Cells.Item(5, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(5, 5).Value = "SVALÖV"
$ws.Cells.Item(5, 7).Value = 5.6
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = 2
$ws.Cells.Item(5, 18).Value = "Ask`r`nSkogsknipprot"
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 61064-2024 artfynd.xlsx", "A 61064-2024")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 61064-2024 karta.png", "A 61064-2024")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 61064-2024 FSC-klagomål.docx", "A 61064-2024")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 61064-2024 FSC-klagomål mail.docx", "A 61064-2024")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 61064-2024 tillsynsbegäran.docx", "A 61064-2024")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 61064-2024 tillsynsbegäran mail.docx", "A 61064-2024")'

# Row 6: A 26855-2022
$ws.Cells.Item(6, 1).Value = "A 26855-2022"
$ws.Cells.Item(6, 2).Value = 44740
$ws.Cells.Item(6, 3).Value = 46081
$ws.Cells.Item(6, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(6, 5).Value = "SVALÖV"
$ws.Cells.Item(6, 7).Value = 5.2
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 2
$ws.Cells.Item(6, 18).Value = "Skogsbräsma`r`nMattlummer"
$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 26855-2022 artfynd.xlsx", "A 26855-2022")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 26855-2022 karta.png", "A 26855-2022")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 26855-2022 FSC-klagomål.docx", "A 26855-2022")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 26855-2022 FSC-klagomål mail.docx", "A 26855-2022")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 26855-2022 tillsynsbegäran.docx", "A 26855-2022")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 26855-2022 tillsynsbegäran mail.docx", "A 26855-2022")'

# Row 7: A 15456-2024
$ws.Cells.Item(7, 1).Value = "A 15456-2024"
$ws.Cells.Item(7, 2).Value = 45401
$ws.Cells.Item(7, 3).Value = 46081
$ws.Cells.Item(7, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(7, 5).Value = "SVALÖV"
$ws.Cells.Item(7, 7).Value = 2.5
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 1
$ws.Cells.Item(7, 16).Value = 1
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Lundticka"
$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 15456-2024 artfynd.xlsx", "A 15456-2024")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 15456-2024 karta.png", "A 15456-2024")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 15456-2024 FSC-klagomål.docx", "A 15456-2024")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 15456-2024 FSC-klagomål mail.docx", "A 15456-2024")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 15456-2024 tillsynsbegäran.docx", "A 15456-2024")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 15456-2024 tillsynsbegäran mail.docx", "A 15456-2024")'

# Row 8: A 15475-2024
$ws.Cells.Item(8, 1).Value = "A 15475-2024"
$ws.Cells.Item(8, 2).Value = 45401
$ws.Cells.Item(8, 3).Value = 46081
$ws.Cells.Item(8, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(8, 5).Value = "SVALÖV"
$ws.Cells.Item(8, 7).Value = 4.7
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Strutbräken"
$ws.Cells.Item(8, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 15475-2024 artfynd.xlsx", "A 15475-2024")'
$ws.Cells.Item(8, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 15475-2024 karta.png", "A 15475-2024")'
$ws.Cells.Item(8, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 15475-2024 FSC-klagomål.docx", "A 15475-2024")'
$ws.Cells.Item(8, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 15475-2024 FSC-klagomål mail.docx", "A 15475-2024")'
$ws.Cells.Item(8, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 15475-2024 tillsynsbegäran.docx", "A 15475-2024")'
$ws.Cells.Item(8, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 15475-2024 tillsynsbegäran mail.docx", "A 15475-2024")'

# Row 9: A 54424-2023
$ws.Cells.Item(9, 1).Value = "A 54424-2023"
$ws.Cells.Item(9, 2).Value = 45233
$ws.Cells.Item(9, 3).Value = 46081
$ws.Cells.Item(9, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(9, 5).Value = "SVALÖV"
$ws.Cells.Item(9, 7).Value = 5.7
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Myskmadra"
$ws.Cells.Item(9, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 54424-2023 artfynd.xlsx", "A 54424-2023")'
$ws.Cells.Item(9, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 54424-2023 karta.png", "A 54424-2023")'
$ws.Cells.Item(9, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 54424-2023 FSC-klagomål.docx", "A 54424-2023")'
$ws.Cells.Item(9, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 54424-2023 FSC-klagomål mail.docx", "A 54424-2023")'
$ws.Cells.Item(9, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 54424-2023 tillsynsbegäran.docx", "A 54424-2023")'
$ws.Cells.Item(9, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 54424-2023 tillsynsbegäran mail.docx", "A 54424-2023")'

# Row 10: A 1577-2024
$ws.Cells.Item(10, 1).Value = "A 1577-2024"
$ws.Cells.Item(10, 2).Value = 45306
$ws.Cells.Item(10, 3).Value = 46081
$ws.Cells.Item(10, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(10, 5).Value = "SVALÖV"
$ws.Cells.Item(10, 7).Value = 21.2
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Större vattensalamander"
$ws.Cells.Item(10, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 1577-2024 artfynd.xlsx", "A 1577-2024")'
$ws.Cells.Item(10, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 1577-2024 karta.png", "A 1577-2024")'
$ws.Cells.Item(10, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 1577-2024 FSC-klagomål.docx", "A 1577-2024")'
$ws.Cells.Item(10, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 1577-2024 FSC-klagomål mail.docx", "A 1577-2024")'
$ws.Cells.Item(10, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 1577-2024 tillsynsbegäran.docx", "A 1577-2024")'
$ws.Cells.Item(10, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 1577-2024 tillsynsbegäran mail.docx", "A 1577-2024")'

# Row 11: A 11170-2023
$ws.Cells.Item(11, 1).Value = "A 11170-2023"
$ws.Cells.Item(11, 2).Value = 44987
$ws.Cells.Item(11, 3).Value = 46081
$ws.Cells.Item(11, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(11, 5).Value = "SVALÖV"
$ws.Cells.Item(11, 7).Value = 0.9
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(11, 16).Value = 1
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Korndådra"
$ws.Cells.Item(11, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 11170-2023 artfynd.xlsx", "A 11170-2023")'
$ws.Cells.Item(11, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 11170-2023 karta.png", "A 11170-2023")'
$ws.Cells.Item(11, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 11170-2023 FSC-klagomål.docx", "A 11170-2023")'
$ws.Cells.Item(11, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 11170-2023 FSC-klagomål mail.docx", "A 11170-2023")'
$ws.Cells.Item(11, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 11170-2023 tillsynsbegäran.docx", "A 11170-2023")'
$ws.Cells.Item(11, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 11170-2023 tillsynsbegäran mail.docx", "A 11170-2023")'

# Row 12: A 54127-2025
$ws.Cells.Item(12, 1).Value = "A 54127-2025"
$ws.Cells.Item(12, 2).Value = 45964
$ws.Cells.Item(12, 3).Value = 46081
$ws.Cells.Item(12, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(12, 5).Value = "SVALÖV"
$ws.Cells.Item(12, 7).Value = 1.3
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Entita"
$ws.Cells.Item(12, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 54127-2025 artfynd.xlsx", "A 54127-2025")'
$ws.Cells.Item(12, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 54127-2025 karta.png", "A 54127-2025")'
$ws.Cells.Item(12, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 54127-2025 FSC-klagomål.docx", "A 54127-2025")'
$ws.Cells.Item(12, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 54127-2025 FSC-klagomål mail.docx", "A 54127-2025")'
$ws.Cells.Item(12, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 54127-2025 tillsynsbegäran.docx", "A 54127-2025")'
$ws.Cells.Item(12, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 54127-2025 tillsynsbegäran mail.docx", "A 54127-2025")'
$ws.Cells.Item(12, 26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/fåglar/A 54127-2025 prioriterade fågelarter.docx", "A 54127-2025")'

# Row 13: A 8169-2024
$ws.Cells.Item(13, 1).Value = "A 8169-2024"
$ws.Cells.Item(13, 2).Value = 45351
$ws.Cells.Item(13, 3).Value = 46081
$ws.Cells.Item(13, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(13, 5).Value = "SVALÖV"
$ws.Cells.Item(13, 7).Value = 3.1
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 1
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 1
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Mjukdån"
$ws.Cells.Item(13, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 8169-2024 artfynd.xlsx", "A 8169-2024")'
$ws.Cells.Item(13, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 8169-2024 karta.png", "A 8169-2024")'
$ws.Cells.Item(13, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 8169-2024 FSC-klagomål.docx", "A 8169-2024")'
$ws.Cells.Item(13, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 8169-2024 FSC-klagomål mail.docx", "A 8169-2024")'
$ws.Cells.Item(13, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 8169-2024 tillsynsbegäran.docx", "A 8169-2024")'
$ws.Cells.Item(13, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 8169-2024 tillsynsbegäran mail.docx", "A 8169-2024")'

# Row 14: A 30657-2022
$ws.Cells.Item(14, 1).Value = "A 30657-2022"
$ws.Cells.Item(14, 2).Value = 44763.513020833336
$ws.Cells.Item(14, 3).Value = 46081
$ws.Cells.Item(14, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(14, 5).Value = "SVALÖV"
$ws.Cells.Item(14, 7).Value = 0.6
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0

# Row 15: A 5937-2022
$ws.Cells.Item(15, 1).Value = "A 5937-2022"
$ws.Cells.Item(15, 2).Value = 44598
$ws.Cells.Item(15, 3).Value = 46081
$ws.Cells.Item(15, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(15, 5).Value = "SVALÖV"
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0

# Row 16: A 57451-2021
$ws.Cells.Item(16, 1).Value = "A 57451-2021"
$ws.Cells.Item(16, 2).Value = 44483.650509259256
$ws.Cells.Item(16, 3).Value = 46081
$ws.Cells.Item(16, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(16, 5).Value = "SVALÖV"
$ws.Cells.Item(16, 7).Value = 0.3
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0

# Row 17: A 46254-2022
$ws.Cells.Item(17, 1).Value = "A 46254-2022"
$ws.Cells.Item(17, 2).Value = 44847.63377314815
$ws.Cells.Item(17, 3).Value = 46081
$ws.Cells.Item(17, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(17, 5).Value = "SVALÖV"
$ws.Cells.Item(17, 7).Value = 0.5
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0

# Row 18: A 15007-2021
$ws.Cells.Item(18, 1).Value = "A 15007-2021"
$ws.Cells.Item(18, 2).Value = 44281
$ws.Cells.Item(18, 3).Value = 46081
$ws.Cells.Item(18, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(18, 5).Value = "SVALÖV"
$ws.Cells.Item(18, 7).Value = 1.7
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 0

# Row 19: A 61495-2021
$ws.Cells.Item(19, 1).Value = "A 61495-2021"
$ws.Cells.Item(19, 2).Value = 44501
$ws.Cells.Item(19, 3).Value = 46081
$ws.Cells.Item(19, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(19, 5).Value = "SVALÖV"
$ws.Cells.Item(19, 7).Value = 1.1
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 0

# Row 20: A 17006-2021
$ws.Cells.Item(20, 1).Value = "A 17006-2021"
$ws.Cells.Item(20, 2).Value = 44295
$ws.Cells.Item(20, 3).Value = 46081
$ws.Cells.Item(20, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(20, 5).Value = "SVALÖV"
$ws.Cells.Item(20, 7).Value = 8.7
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 0

# Row 21: A 42617-2022
$ws.Cells.Item(21, 1).Value = "A 42617-2022"
$ws.Cells.Item(21, 2).Value = 44831.858819444446
$ws.Cells.Item(21, 3).Value = 46081
$ws.Cells.Item(21, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(21, 5).Value = "SVALÖV"
$ws.Cells.Item(21, 7).Value = 0.6
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 0

# Row 22: A 19159-2021
$ws.Cells.Item(22, 1).Value = "A 19159-2021"
$ws.Cells.Item(22, 2).Value = 44308
$ws.Cells.Item(22, 3).Value = 46081
$ws.Cells.Item(22, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(22, 5).Value = "SVALÖV"
$ws.Cells.Item(22, 7).Value = 0.8
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 0
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 0

# Row 23: A 47717-2022
$ws.Cells.Item(23, 1).Value = "A 47717-2022"
$ws.Cells.Item(23, 2).Value = 44854.605844907404
$ws.Cells.Item(23, 3).Value = 46081
$ws.Cells.Item(23, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(23, 5).Value = "SVALÖV"
$ws.Cells.Item(23, 7).Value = 0.6
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(23, 15).Value = 0
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 0

# Row 24: A 13245-2024
$ws.Cells.Item(24, 1).Value = "A 13245-2024"
$ws.Cells.Item(24, 2).Value = 45386
$ws.Cells.Item(24, 3).Value = 46081
$ws.Cells.Item(24, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(24, 5).Value = "SVALÖV"
$ws.Cells.Item(24, 7).Value = 13
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0

# Row 25: A 34859-2021
$ws.Cells.Item(25, 1).Value = "A 34859-2021"
$ws.Cells.Item(25, 2).Value = 44382
$ws.Cells.Item(25, 3).Value = 46081
$ws.Cells.Item(25, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(25, 5).Value = "SVALÖV"
$ws.Cells.Item(25, 7).Value = 1.9
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 0

# Row 26: A 21774-2024
$ws.Cells.Item(26, 1).Value = "A 21774-2024"
$ws.Cells.Item(26, 2).Value = 45442
$ws.Cells.Item(26, 3).Value = 46081
$ws.Cells.Item(26, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(26, 5).Value = "SVALÖV"
$ws.Cells.Item(26, 7).Value = 0.6
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = 0
$ws.Cells.Item(26, 14).Value = 0
$ws.Cells.Item(26, 15).Value = 0
$ws.Cells.Item(26, 16).Value = 0
$ws.Cells.Item(26, 17).Value = 0

# Row 27: A 34810-2023
$ws.Cells.Item(27, 1).Value = "A 34810-2023"
$ws.Cells.Item(27, 2).Value = 45141
$ws.Cells.Item(27, 3).Value = 46081
$ws.Cells.Item(27, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(27, 5).Value = "SVALÖV"
$ws.Cells.Item(27, 7).Value = 7
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 0
$ws.Cells.Item(27, 14).Value = 0
$ws.Cells.Item(27, 15).Value = 0
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(27, 17).Value = 0

# Row 28: A 22794-2024
$ws.Cells.Item(28, 1).Value = "A 22794-2024"
$ws.Cells.Item(28, 2).Value = 45448.49752314815
$ws.Cells.Item(28, 3).Value = 46081
$ws.Cells.Item(28, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(28, 5).Value = "SVALÖV"
$ws.Cells.Item(28, 7).Value = 3.8
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(28, 14).Value = 0
$ws.Cells.Item(28, 15).Value = 0
$ws.Cells.Item(28, 16).Value = 0
$ws.Cells.Item(28, 17).Value = 0

# Row 29: A 35587-2024
$ws.Cells.Item(29, 1).Value = "A 35587-2024"
$ws.Cells.Item(29, 2).Value = 45531
$ws.Cells.Item(29, 3).Value = 46081
$ws.Cells.Item(29, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(29, 5).Value = "SVALÖV"
$ws.Cells.Item(29, 7).Value = 1.6
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

# Row 30: A 38394-2025
$ws.Cells.Item(30, 1).Value = "A 38394-2025"
$ws.Cells.Item(30, 2).Value = 45883
$ws.Cells.Item(30, 3).Value = 46081
$ws.Cells.Item(30, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(30, 5).Value = "SVALÖV"
$ws.Cells.Item(30, 7).Value = 0.5
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = 0
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0

# Row 31: A 16667-2023
$ws.Cells.Item(31, 1).Value = "A 16667-2023"
$ws.Cells.Item(31, 2).Value = 45030
$ws.Cells.Item(31, 3).Value = 46081
$ws.Cells.Item(31, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(31, 5).Value = "SVALÖV"
$ws.Cells.Item(31, 7).Value = 3.8
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(31, 14).Value = 0
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(31, 17).Value = 0

# Row 32: A 26949-2025
$ws.Cells.Item(32, 1).Value = "A 26949-2025"
$ws.Cells.Item(32, 2).Value = 45811.43701388889
$ws.Cells.Item(32, 3).Value = 46081
$ws.Cells.Item(32, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(32, 5).Value = "SVALÖV"
$ws.Cells.Item(32, 7).Value = 4.2
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = 0
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(32, 17).Value = 0

# Row 33: A 47696-2024
$ws.Cells.Item(33, 1).Value = "A 47696-2024"
$ws.Cells.Item(33, 2).Value = 45588.458553240744
$ws.Cells.Item(33, 3).Value = 46081
$ws.Cells.Item(33, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(33, 5).Value = "SVALÖV"
$ws.Cells.Item(33, 7).Value = 3.3
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 0
$ws.Cells.Item(33, 14).Value = 0
$ws.Cells.Item(33, 15).Value = 0
$ws.Cells.Item(33, 16).Value = 0
$ws.Cells.Item(33, 17).Value = 0

# Row 34: A 8332-2025
$ws.Cells.Item(34, 1).Value = "A 8332-2025"
$ws.Cells.Item(34, 2).Value = 45708.648194444446
$ws.Cells.Item(34, 3).Value = 46081
$ws.Cells.Item(34, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(34, 5).Value = "SVALÖV"
$ws.Cells.Item(34, 7).Value = 0.5
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(34, 14).Value = 0
$ws.Cells.Item(34, 15).Value = 0
$ws.Cells.Item(34, 16).Value = 0
$ws.Cells.Item(34, 17).Value = 0

# Row 35: A 14543-2024
$ws.Cells.Item(35, 1).Value = "A 14543-2024"
$ws.Cells.Item(35, 2).Value = 45394
$ws.Cells.Item(35, 3).Value = 46081
$ws.Cells.Item(35, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(35, 5).Value = "SVALÖV"
$ws.Cells.Item(35, 7).Value = 0.7
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = 0
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(35, 17).Value = 0

# Row 36: A 42706-2025
$ws.Cells.Item(36, 1).Value = "A 42706-2025"
$ws.Cells.Item(36, 2).Value = 45908.37222222222
$ws.Cells.Item(36, 3).Value = 46081
$ws.Cells.Item(36, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(36, 5).Value = "SVALÖV"
$ws.Cells.Item(36, 7).Value = 0.5
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(36, 17).Value = 0

# Row 37: A 22063-2025
$ws.Cells.Item(37, 1).Value = "A 22063-2025"
$ws.Cells.Item(37, 2).Value = 45785.37700231482
$ws.Cells.Item(37, 3).Value = 46081
$ws.Cells.Item(37, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(37, 5).Value = "SVALÖV"
$ws.Cells.Item(37, 7).Value = 2
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = 0
$ws.Cells.Item(37, 14).Value = 0
$ws.Cells.Item(37, 15).Value = 0
$ws.Cells.Item(37, 16).Value = 0
$ws.Cells.Item(37, 17).Value = 0

# Row 38: A 42698-2025
$ws.Cells.Item(38, 1).Value = "A 42698-2025"
$ws.Cells.Item(38, 2).Value = 45908.3650462963
$ws.Cells.Item(38, 3).Value = 46081
$ws.Cells.Item(38, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(38, 5).Value = "SVALÖV"
$ws.Cells.Item(38, 7).Value = 3.1
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 0
$ws.Cells.Item(38, 14).Value = 0
$ws.Cells.Item(38, 15).Value = 0
$ws.Cells.Item(38, 16).Value = 0
$ws.Cells.Item(38, 17).Value = 0

# Row 39: A 20170-2024
$ws.Cells.Item(39, 1).Value = "A 20170-2024"
$ws.Cells.Item(39, 2).Value = 45434
$ws.Cells.Item(39, 3).Value = 46081
$ws.Cells.Item(39, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(39, 5).Value = "SVALÖV"
$ws.Cells.Item(39, 7).Value = 0.6
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = 0
$ws.Cells.Item(39, 14).Value = 0
$ws.Cells.Item(39, 15).Value = 0
$ws.Cells.Item(39, 16).Value = 0
$ws.Cells.Item(39, 17).Value = 0

# Row 40: A 62403-2022
$ws.Cells.Item(40, 1).Value = "A 62403-2022"
$ws.Cells.Item(40, 2).Value = 44923
$ws.Cells.Item(40, 3).Value = 46081
$ws.Cells.Item(40, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(40, 5).Value = "SVALÖV"
$ws.Cells.Item(40, 7).Value = 2.2
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = 0
$ws.Cells.Item(40, 15).Value = 0
$ws.Cells.Item(40, 16).Value = 0
$ws.Cells.Item(40, 17).Value = 0

# Row 41: A 43425-2025
$ws.Cells.Item(41, 1).Value = "A 43425-2025"
$ws.Cells.Item(41, 2).Value = 45911.418217592596
$ws.Cells.Item(41, 3).Value = 46081
$ws.Cells.Item(41, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(41, 5).Value = "SVALÖV"
$ws.Cells.Item(41, 7).Value = 2.1
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = 0
$ws.Cells.Item(41, 14).Value = 0
$ws.Cells.Item(41, 15).Value = 0
$ws.Cells.Item(41, 16).Value = 0
$ws.Cells.Item(41, 17).Value = 0

# Row 42: A 19626-2021
$ws.Cells.Item(42, 1).Value = "A 19626-2021"
$ws.Cells.Item(42, 2).Value = 44309
$ws.Cells.Item(42, 3).Value = 46081
$ws.Cells.Item(42, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(42, 5).Value = "SVALÖV"
$ws.Cells.Item(42, 7).Value = 17.8
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(42, 17).Value = 0

# Row 43: A 17980-2025
$ws.Cells.Item(43, 1).Value = "A 17980-2025"
$ws.Cells.Item(43, 2).Value = 45761.36854166666
$ws.Cells.Item(43, 3).Value = 46081
$ws.Cells.Item(43, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(43, 5).Value = "SVALÖV"
$ws.Cells.Item(43, 7).Value = 3.4
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(43, 14).Value = 0
$ws.Cells.Item(43, 15).Value = 0
$ws.Cells.Item(43, 16).Value = 0
$ws.Cells.Item(43, 17).Value = 0

# Row 44: A 19282-2025
$ws.Cells.Item(44, 1).Value = "A 19282-2025"
$ws.Cells.Item(44, 2).Value = 45769
$ws.Cells.Item(44, 3).Value = 46081
$ws.Cells.Item(44, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(44, 5).Value = "SVALÖV"
$ws.Cells.Item(44, 7).Value = 0.8
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(44, 17).Value = 0

# Row 45: A 15157-2023
$ws.Cells.Item(45, 1).Value = "A 15157-2023"
$ws.Cells.Item(45, 2).Value = 45016
$ws.Cells.Item(45, 3).Value = 46081
$ws.Cells.Item(45, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(45, 5).Value = "SVALÖV"
$ws.Cells.Item(45, 7).Value = 2
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(45, 17).Value = 0

# Row 46: A 11491-2023
$ws.Cells.Item(46, 1).Value = "A 11491-2023"
$ws.Cells.Item(46, 2).Value = 44991
$ws.Cells.Item(46, 3).Value = 46081
$ws.Cells.Item(46, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(46, 5).Value = "SVALÖV"
$ws.Cells.Item(46, 7).Value = 2.8
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 0
$ws.Cells.Item(46, 14).Value = 0
$ws.Cells.Item(46, 15).Value = 0
$ws.Cells.Item(46, 16).Value = 0
$ws.Cells.Item(46, 17).Value = 0

# Row 47: A 29296-2021
$ws.Cells.Item(47, 1).Value = "A 29296-2021"
$ws.Cells.Item(47, 2).Value = 44361
$ws.Cells.Item(47, 3).Value = 46081
$ws.Cells.Item(47, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(47, 5).Value = "SVALÖV"
$ws.Cells.Item(47, 7).Value = 0.7
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0
$ws.Cells.Item(47, 14).Value = 0
$ws.Cells.Item(47, 15).Value = 0
$ws.Cells.Item(47, 16).Value = 0
$ws.Cells.Item(47, 17).Value = 0

# Row 48: A 10209-2023
$ws.Cells.Item(48, 1).Value = "A 10209-2023"
$ws.Cells.Item(48, 2).Value = 44986
$ws.Cells.Item(48, 3).Value = 46081
$ws.Cells.Item(48, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(48, 5).Value = "SVALÖV"
$ws.Cells.Item(48, 7).Value = 1.2
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 0
$ws.Cells.Item(48, 14).Value = 0
$ws.Cells.Item(48, 15).Value = 0
$ws.Cells.Item(48, 16).Value = 0
$ws.Cells.Item(48, 17).Value = 0

# Row 49: A 47139-2025
$ws.Cells.Item(49, 1).Value = "A 47139-2025"
$ws.Cells.Item(49, 2).Value = 45929
$ws.Cells.Item(49, 3).Value = 46081
$ws.Cells.Item(49, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(49, 5).Value = "SVALÖV"
$ws.Cells.Item(49, 7).Value = 1.1
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = 0
$ws.Cells.Item(49, 14).Value = 0
$ws.Cells.Item(49, 15).Value = 0
$ws.Cells.Item(49, 16).Value = 0
$ws.Cells.Item(49, 17).Value = 0

# Row 50: A 37907-2025
$ws.Cells.Item(50, 1).Value = "A 37907-2025"
$ws.Cells.Item(50, 2).Value = 45880
$ws.Cells.Item(50, 3).Value = 46081
$ws.Cells.Item(50, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(50, 5).Value = "SVALÖV"
$ws.Cells.Item(50, 6).Value = "Kyrkan"
$ws.Cells.Item(50, 7).Value = 5.7
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(50, 14).Value = 0
$ws.Cells.Item(50, 15).Value = 0
$ws.Cells.Item(50, 16).Value = 0
$ws.Cells.Item(50, 17).Value = 0

# Row 51: A 44107-2023
$ws.Cells.Item(51, 1).Value = "A 44107-2023"
$ws.Cells.Item(51, 2).Value = 45188.43670138889
$ws.Cells.Item(51, 3).Value = 46081
$ws.Cells.Item(51, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(51, 5).Value = "SVALÖV"
$ws.Cells.Item(51, 7).Value = 2.3
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(51, 14).Value = 0
$ws.Cells.Item(51, 15).Value = 0
$ws.Cells.Item(51, 16).Value = 0
$ws.Cells.Item(51, 17).Value = 0

# Row 52: A 40470-2025
$ws.Cells.Item(52, 1).Value = "A 40470-2025"
$ws.Cells.Item(52, 2).Value = 45895
$ws.Cells.Item(52, 3).Value = 46081
$ws.Cells.Item(52, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(52, 5).Value = "SVALÖV"
$ws.Cells.Item(52, 7).Value = 0.1
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 13).Value = 0
$ws.Cells.Item(52, 14).Value = 0
$ws.Cells.Item(52, 15).Value = 0
$ws.Cells.Item(52, 16).Value = 0
$ws.Cells.Item(52, 17).Value = 0

# Row 53: A 20235-2025
$ws.Cells.Item(53, 1).Value = "A 20235-2025"
$ws.Cells.Item(53, 2).Value = 45772.67386574074
$ws.Cells.Item(53, 3).Value = 46081
$ws.Cells.Item(53, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(53, 5).Value = "SVALÖV"
$ws.Cells.Item(53, 7).Value = 0.9
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(53, 14).Value = 0
$ws.Cells.Item(53, 15).Value = 0
$ws.Cells.Item(53, 16).Value = 0
$ws.Cells.Item(53, 17).Value = 0

# Row 54: A 22262-2023
$ws.Cells.Item(54, 1).Value = "A 22262-2023"
$ws.Cells.Item(54, 2).Value = 45070
$ws.Cells.Item(54, 3).Value = 46081
$ws.Cells.Item(54, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(54, 5).Value = "SVALÖV"
$ws.Cells.Item(54, 7).Value = 1.1
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = 0
$ws.Cells.Item(54, 14).Value = 0
$ws.Cells.Item(54, 15).Value = 0
$ws.Cells.Item(54, 16).Value = 0
$ws.Cells.Item(54, 17).Value = 0

# Row 55: A 48279-2021
$ws.Cells.Item(55, 1).Value = "A 48279-2021"
$ws.Cells.Item(55, 2).Value = 44449
$ws.Cells.Item(55, 3).Value = 46081
$ws.Cells.Item(55, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(55, 5).Value = "SVALÖV"
$ws.Cells.Item(55, 7).Value = 1.3
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = 0
$ws.Cells.Item(55, 14).Value = 0
$ws.Cells.Item(55, 15).Value = 0
$ws.Cells.Item(55, 16).Value = 0
$ws.Cells.Item(55, 17).Value = 0

# Row 56: A 49032-2022
$ws.Cells.Item(56, 1).Value = "A 49032-2022"
$ws.Cells.Item(56, 2).Value = 44860.440833333334
$ws.Cells.Item(56, 3).Value = 46081
$ws.Cells.Item(56, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(56, 5).Value = "SVALÖV"
$ws.Cells.Item(56, 7).Value = 3.7
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = 0
$ws.Cells.Item(56, 14).Value = 0
$ws.Cells.Item(56, 15).Value = 0
$ws.Cells.Item(56, 16).Value = 0
$ws.Cells.Item(56, 17).Value = 0

# Row 57: A 18121-2025
$ws.Cells.Item(57, 1).Value = "A 18121-2025"
$ws.Cells.Item(57, 2).Value = 45761
$ws.Cells.Item(57, 3).Value = 46081
$ws.Cells.Item(57, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(57, 5).Value = "SVALÖV"
$ws.Cells.Item(57, 7).Value = 1.7
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = 0
$ws.Cells.Item(57, 14).Value = 0
$ws.Cells.Item(57, 15).Value = 0
$ws.Cells.Item(57, 16).Value = 0
$ws.Cells.Item(57, 17).Value = 0

# Row 58: A 43294-2025
$ws.Cells.Item(58, 1).Value = "A 43294-2025"
$ws.Cells.Item(58, 2).Value = 45910
$ws.Cells.Item(58, 3).Value = 46081
$ws.Cells.Item(58, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(58, 5).Value = "SVALÖV"
$ws.Cells.Item(58, 7).Value = 2.6
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = 0
$ws.Cells.Item(58, 14).Value = 0
$ws.Cells.Item(58, 15).Value = 0
$ws.Cells.Item(58, 16).Value = 0
$ws.Cells.Item(58, 17).Value = 0

# Row 59: A 6335-2022
$ws.Cells.Item(59, 1).Value = "A 6335-2022"
$ws.Cells.Item(59, 2).Value = 44600
$ws.Cells.Item(59, 3).Value = 46081
$ws.Cells.Item(59, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(59, 5).Value = "SVALÖV"
$ws.Cells.Item(59, 7).Value = 4
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 13).Value = 0
$ws.Cells.Item(59, 14).Value = 0
$ws.Cells.Item(59, 15).Value = 0
$ws.Cells.Item(59, 16).Value = 0
$ws.Cells.Item(59, 17).Value = 0

# Row 60: A 20125-2024
$ws.Cells.Item(60, 1).Value = "A 20125-2024"
$ws.Cells.Item(60, 2).Value = 45434.55394675926
$ws.Cells.Item(60, 3).Value = 46081
$ws.Cells.Item(60, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(60, 5).Value = "SVALÖV"
$ws.Cells.Item(60, 7).Value = 5.7
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = 0
$ws.Cells.Item(60, 14).Value = 0
$ws.Cells.Item(60, 15).Value = 0
$ws.Cells.Item(60, 16).Value = 0
$ws.Cells.Item(60, 17).Value = 0

# Row 61: A 10196-2025
$ws.Cells.Item(61, 1).Value = "A 10196-2025"
$ws.Cells.Item(61, 2).Value = 45719
$ws.Cells.Item(61, 3).Value = 46081
$ws.Cells.Item(61, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(61, 5).Value = "SVALÖV"
$ws.Cells.Item(61, 7).Value = 0.5
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = 0
$ws.Cells.Item(61, 14).Value = 0
$ws.Cells.Item(61, 15).Value = 0
$ws.Cells.Item(61, 16).Value = 0
$ws.Cells.Item(61, 17).Value = 0

# Row 62: A 58970-2025
$ws.Cells.Item(62, 1).Value = "A 58970-2025"
$ws.Cells.Item(62, 2).Value = 45987.61269675926
$ws.Cells.Item(62, 3).Value = 46081
$ws.Cells.Item(62, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(62, 5).Value = "SVALÖV"
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = 0
$ws.Cells.Item(62, 14).Value = 0
$ws.Cells.Item(62, 15).Value = 0
$ws.Cells.Item(62, 16).Value = 0
$ws.Cells.Item(62, 17).Value = 0

# Row 63: A 59230-2025
$ws.Cells.Item(63, 1).Value = "A 59230-2025"
$ws.Cells.Item(63, 2).Value = 45988.625023148146
$ws.Cells.Item(63, 3).Value = 46081
$ws.Cells.Item(63, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(63, 5).Value = "SVALÖV"
$ws.Cells.Item(63, 7).Value = 0.6
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(63, 14).Value = 0
$ws.Cells.Item(63, 15).Value = 0
$ws.Cells.Item(63, 16).Value = 0
$ws.Cells.Item(63, 17).Value = 0

# Row 64: A 44112-2023
$ws.Cells.Item(64, 1).Value = "A 44112-2023"
$ws.Cells.Item(64, 2).Value = 45188.441828703704
$ws.Cells.Item(64, 3).Value = 46081
$ws.Cells.Item(64, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(64, 5).Value = "SVALÖV"
$ws.Cells.Item(64, 7).Value = 0.5
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(64, 14).Value = 0
$ws.Cells.Item(64, 15).Value = 0
$ws.Cells.Item(64, 16).Value = 0
$ws.Cells.Item(64, 17).Value = 0

# Row 65: A 43466-2025
$ws.Cells.Item(65, 1).Value = "A 43466-2025"
$ws.Cells.Item(65, 2).Value = 45911
$ws.Cells.Item(65, 3).Value = 46081
$ws.Cells.Item(65, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(65, 5).Value = "SVALÖV"
$ws.Cells.Item(65, 7).Value = 4.5
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(65, 14).Value = 0
$ws.Cells.Item(65, 15).Value = 0
$ws.Cells.Item(65, 16).Value = 0
$ws.Cells.Item(65, 17).Value = 0

# Row 66: A 35911-2023
$ws.Cells.Item(66, 1).Value = "A 35911-2023"
$ws.Cells.Item(66, 2).Value = 45148.64078703704
$ws.Cells.Item(66, 3).Value = 46081
$ws.Cells.Item(66, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(66, 5).Value = "SVALÖV"
$ws.Cells.Item(66, 7).Value = 0.9
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = 0
$ws.Cells.Item(66, 14).Value = 0
$ws.Cells.Item(66, 15).Value = 0
$ws.Cells.Item(66, 16).Value = 0
$ws.Cells.Item(66, 17).Value = 0

# Row 67: A 11672-2024
$ws.Cells.Item(67, 1).Value = "A 11672-2024"
$ws.Cells.Item(67, 2).Value = 45373
$ws.Cells.Item(67, 3).Value = 46081
$ws.Cells.Item(67, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(67, 5).Value = "SVALÖV"
$ws.Cells.Item(67, 7).Value = 13.9
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = 0
$ws.Cells.Item(67, 14).Value = 0
$ws.Cells.Item(67, 15).Value = 0
$ws.Cells.Item(67, 16).Value = 0
$ws.Cells.Item(67, 17).Value = 0

# Row 68: A 42588-2023
$ws.Cells.Item(68, 1).Value = "A 42588-2023"
$ws.Cells.Item(68, 2).Value = 45176
$ws.Cells.Item(68, 3).Value = 46081
$ws.Cells.Item(68, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(68, 5).Value = "SVALÖV"
$ws.Cells.Item(68, 7).Value = 1.1
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(68, 14).Value = 0
$ws.Cells.Item(68, 15).Value = 0
$ws.Cells.Item(68, 16).Value = 0
$ws.Cells.Item(68, 17).Value = 0

# Row 69: A 1901-2025
$ws.Cells.Item(69, 1).Value = "A 1901-2025"
$ws.Cells.Item(69, 2).Value = 45671
$ws.Cells.Item(69, 3).Value = 46081
$ws.Cells.Item(69, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(69, 5).Value = "SVALÖV"
$ws.Cells.Item(69, 7).Value = 4.2
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(69, 14).Value = 0
$ws.Cells.Item(69, 15).Value = 0
$ws.Cells.Item(69, 16).Value = 0
$ws.Cells.Item(69, 17).Value = 0

# Row 70: A 1920-2025
$ws.Cells.Item(70, 1).Value = "A 1920-2025"
$ws.Cells.Item(70, 2).Value = 45671
$ws.Cells.Item(70, 3).Value = 46081
$ws.Cells.Item(70, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(70, 5).Value = "SVALÖV"
$ws.Cells.Item(70, 7).Value = 2.3
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(70, 14).Value = 0
$ws.Cells.Item(70, 15).Value = 0
$ws.Cells.Item(70, 16).Value = 0
$ws.Cells.Item(70, 17).Value = 0

# Row 71: A 43736-2023
$ws.Cells.Item(71, 1).Value = "A 43736-2023"
$ws.Cells.Item(71, 2).Value = 45182
$ws.Cells.Item(71, 3).Value = 46081
$ws.Cells.Item(71, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(71, 5).Value = "SVALÖV"
$ws.Cells.Item(71, 7).Value = 1.2
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(71, 14).Value = 0
$ws.Cells.Item(71, 15).Value = 0
$ws.Cells.Item(71, 16).Value = 0
$ws.Cells.Item(71, 17).Value = 0

# Row 72: A 43004-2025
$ws.Cells.Item(72, 1).Value = "A 43004-2025"
$ws.Cells.Item(72, 2).Value = 45909
$ws.Cells.Item(72, 3).Value = 46081
$ws.Cells.Item(72, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(72, 5).Value = "SVALÖV"
$ws.Cells.Item(72, 7).Value = 1.6
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = 0
$ws.Cells.Item(72, 14).Value = 0
$ws.Cells.Item(72, 15).Value = 0
$ws.Cells.Item(72, 16).Value = 0
$ws.Cells.Item(72, 17).Value = 0

# Row 73: A 52664-2023
$ws.Cells.Item(73, 1).Value = "A 52664-2023"
$ws.Cells.Item(73, 2).Value = 45225
$ws.Cells.Item(73, 3).Value = 46081
$ws.Cells.Item(73, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(73, 5).Value = "SVALÖV"
$ws.Cells.Item(73, 7).Value = 6.4
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(73, 14).Value = 0
$ws.Cells.Item(73, 15).Value = 0
$ws.Cells.Item(73, 16).Value = 0
$ws.Cells.Item(73, 17).Value = 0

# Row 74: A 9716-2022
$ws.Cells.Item(74, 1).Value = "A 9716-2022"
$ws.Cells.Item(74, 2).Value = 44617
$ws.Cells.Item(74, 3).Value = 46081
$ws.Cells.Item(74, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(74, 5).Value = "SVALÖV"
$ws.Cells.Item(74, 7).Value = 0.7
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(74, 14).Value = 0
$ws.Cells.Item(74, 15).Value = 0
$ws.Cells.Item(74, 16).Value = 0
$ws.Cells.Item(74, 17).Value = 0

# Row 75: A 61045-2024
$ws.Cells.Item(75, 1).Value = "A 61045-2024"
$ws.Cells.Item(75, 2).Value = 45645.47678240741
$ws.Cells.Item(75, 3).Value = 46081
$ws.Cells.Item(75, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(75, 5).Value = "SVALÖV"
$ws.Cells.Item(75, 7).Value = 5.7
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(75, 14).Value = 0
$ws.Cells.Item(75, 15).Value = 0
$ws.Cells.Item(75, 16).Value = 0
$ws.Cells.Item(75, 17).Value = 0

# Row 76: A 37354-2022
$ws.Cells.Item(76, 1).Value = "A 37354-2022"
$ws.Cells.Item(76, 2).Value = 44806
$ws.Cells.Item(76, 3).Value = 46081
$ws.Cells.Item(76, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(76, 5).Value = "SVALÖV"
$ws.Cells.Item(76, 7).Value = 1.4
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = 0
$ws.Cells.Item(76, 14).Value = 0
$ws.Cells.Item(76, 15).Value = 0
$ws.Cells.Item(76, 16).Value = 0
$ws.Cells.Item(76, 17).Value = 0

# Row 77: A 5190-2024
$ws.Cells.Item(77, 1).Value = "A 5190-2024"
$ws.Cells.Item(77, 2).Value = 45330
$ws.Cells.Item(77, 3).Value = 46081
$ws.Cells.Item(77, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(77, 5).Value = "SVALÖV"
$ws.Cells.Item(77, 7).Value = 2.1
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(77, 14).Value = 0
$ws.Cells.Item(77, 15).Value = 0
$ws.Cells.Item(77, 16).Value = 0
$ws.Cells.Item(77, 17).Value = 0

# Row 78: A 62201-2025
$ws.Cells.Item(78, 1).Value = "A 62201-2025"
$ws.Cells.Item(78, 2).Value = 46006.55501157408
$ws.Cells.Item(78, 3).Value = 46081
$ws.Cells.Item(78, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(78, 5).Value = "SVALÖV"
$ws.Cells.Item(78, 7).Value = 22.1
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).Value = 0
$ws.Cells.Item(78, 14).Value = 0
$ws.Cells.Item(78, 15).Value = 0
$ws.Cells.Item(78, 16).Value = 0
$ws.Cells.Item(78, 17).Value = 0

# Row 79: A 17731-2025
$ws.Cells.Item(79, 1).Value = "A 17731-2025"
$ws.Cells.Item(79, 2).Value = 45758
$ws.Cells.Item(79, 3).Value = 46081
$ws.Cells.Item(79, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(79, 5).Value = "SVALÖV"
$ws.Cells.Item(79, 7).Value = 0.8
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = 0
$ws.Cells.Item(79, 14).Value = 0
$ws.Cells.Item(79, 15).Value = 0
$ws.Cells.Item(79, 16).Value = 0
$ws.Cells.Item(79, 17).Value = 0

# Row 80: A 3430-2026
$ws.Cells.Item(80, 1).Value = "A 3430-2026"
$ws.Cells.Item(80, 2).Value = 46042
$ws.Cells.Item(80, 3).Value = 46081
$ws.Cells.Item(80, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(80, 5).Value = "SVALÖV"
$ws.Cells.Item(80, 7).Value = 5.7
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = 0
$ws.Cells.Item(80, 14).Value = 0
$ws.Cells.Item(80, 15).Value = 0
$ws.Cells.Item(80, 16).Value = 0
$ws.Cells.Item(80, 17).Value = 0

# Row 81: A 5068-2026
$ws.Cells.Item(81, 1).Value = "A 5068-2026"
$ws.Cells.Item(81, 2).Value = 46049.50232638889
$ws.Cells.Item(81, 3).Value = 46081
$ws.Cells.Item(81, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(81, 5).Value = "SVALÖV"
$ws.Cells.Item(81, 7).Value = 3.5
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = 0
$ws.Cells.Item(81, 14).Value = 0
$ws.Cells.Item(81, 15).Value = 0
$ws.Cells.Item(81, 16).Value = 0
$ws.Cells.Item(81, 17).Value = 0

# Row 82: A 9395-2026
$ws.Cells.Item(82, 1).Value = "A 9395-2026"
$ws.Cells.Item(82, 2).Value = 46065
$ws.Cells.Item(82, 3).Value = 46081
$ws.Cells.Item(82, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(82, 5).Value = "SVALÖV"
$ws.Cells.Item(82, 7).Value = 3.7
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = 0
$ws.Cells.Item(82, 14).Value = 0
$ws.Cells.Item(82, 15).Value = 0
$ws.Cells.Item(82, 16).Value = 0
$ws.Cells.Item(82, 17).Value = 0

# Row 83: A 6365-2026
$ws.Cells.Item(83, 1).Value = "A 6365-2026"
$ws.Cells.Item(83, 2).Value = 46051
$ws.Cells.Item(83, 3).Value = 46081
$ws.Cells.Item(83, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(83, 5).Value = "SVALÖV"
$ws.Cells.Item(83, 7).Value = 2.9
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = 0
$ws.Cells.Item(83, 14).Value = 0
$ws.Cells.Item(83, 15).Value = 0
$ws.Cells.Item(83, 16).Value = 0
$ws.Cells.Item(83, 17).Value = 0

# Row 84: A 22792-2024
$ws.Cells.Item(84, 1).Value = "A 22792-2024"
$ws.Cells.Item(84, 2).Value = 45448.48983796296
$ws.Cells.Item(84, 3).Value = 46081
$ws.Cells.Item(84, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(84, 5).Value = "SVALÖV"
$ws.Cells.Item(84, 7).Value = 0.5
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(84, 14).Value = 0
$ws.Cells.Item(84, 15).Value = 0
$ws.Cells.Item(84, 16).Value = 0
$ws.Cells.Item(84, 17).Value = 0

# Row 85: A 59266-2025
$ws.Cells.Item(85, 1).Value = "A 59266-2025"
$ws.Cells.Item(85, 2).Value = 45988
$ws.Cells.Item(85, 3).Value = 46081
$ws.Cells.Item(85, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(85, 5).Value = "SVALÖV"
$ws.Cells.Item(85, 7).Value = 1.8
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(85, 14).Value = 0
$ws.Cells.Item(85, 15).Value = 0
$ws.Cells.Item(85, 16).Value = 0
$ws.Cells.Item(85, 17).Value = 0

# Row 86: A 4964-2022
$ws.Cells.Item(86, 1).Value = "A 4964-2022"
$ws.Cells.Item(86, 2).Value = 44593
$ws.Cells.Item(86, 3).Value = 46081
$ws.Cells.Item(86, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(86, 5).Value = "SVALÖV"
$ws.Cells.Item(86, 7).Value = 1.5
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(86, 14).Value = 0
$ws.Cells.Item(86, 15).Value = 0
$ws.Cells.Item(86, 16).Value = 0
$ws.Cells.Item(86, 17).Value = 0

# Row 87: A 11177-2023
$ws.Cells.Item(87, 1).Value = "A 11177-2023"
$ws.Cells.Item(87, 2).Value = 44987
$ws.Cells.Item(87, 3).Value = 46081
$ws.Cells.Item(87, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(87, 5).Value = "SVALÖV"
$ws.Cells.Item(87, 7).Value = 1.2
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(87, 14).Value = 0
$ws.Cells.Item(87, 15).Value = 0
$ws.Cells.Item(87, 16).Value = 0
$ws.Cells.Item(87, 17).Value = 0

# Row 88: A 16215-2022
$ws.Cells.Item(88, 1).Value = "A 16215-2022"
$ws.Cells.Item(88, 2).Value = 44670
$ws.Cells.Item(88, 3).Value = 46081
$ws.Cells.Item(88, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(88, 5).Value = "SVALÖV"
$ws.Cells.Item(88, 7).Value = 2.9
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(88, 14).Value = 0
$ws.Cells.Item(88, 15).Value = 0
$ws.Cells.Item(88, 16).Value = 0
$ws.Cells.Item(88, 17).Value = 0

# Row 89: A 61620-2022
$ws.Cells.Item(89, 1).Value = "A 61620-2022"
$ws.Cells.Item(89, 2).Value = 44910
$ws.Cells.Item(89, 3).Value = 46081
$ws.Cells.Item(89, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(89, 5).Value = "SVALÖV"
$ws.Cells.Item(89, 7).Value = 1.5
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 0
$ws.Cells.Item(89, 14).Value = 0
$ws.Cells.Item(89, 15).Value = 0
$ws.Cells.Item(89, 16).Value = 0
$ws.Cells.Item(89, 17).Value = 0

# Row 90: A 48262-2021
$ws.Cells.Item(90, 1).Value = "A 48262-2021"
$ws.Cells.Item(90, 2).Value = 44449.66324074074
$ws.Cells.Item(90, 3).Value = 46081
$ws.Cells.Item(90, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(90, 5).Value = "SVALÖV"
$ws.Cells.Item(90, 7).Value = 0.6
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = 0
$ws.Cells.Item(90, 14).Value = 0
$ws.Cells.Item(90, 15).Value = 0
$ws.Cells.Item(90, 16).Value = 0
$ws.Cells.Item(90, 17).Value = 0

# Row 91: A 8924-2023
$ws.Cells.Item(91, 1).Value = "A 8924-2023"
$ws.Cells.Item(91, 2).Value = 44974
$ws.Cells.Item(91, 3).Value = 46081
$ws.Cells.Item(91, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(91, 5).Value = "SVALÖV"
$ws.Cells.Item(91, 7).Value = 6.8
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = 0
$ws.Cells.Item(91, 14).Value = 0
$ws.Cells.Item(91, 15).Value = 0
$ws.Cells.Item(91, 16).Value = 0
$ws.Cells.Item(91, 17).Value = 0

# Row 92: A 497-2023
$ws.Cells.Item(92, 1).Value = "A 497-2023"
$ws.Cells.Item(92, 2).Value = 44930.33540509259
$ws.Cells.Item(92, 3).Value = 46081
$ws.Cells.Item(92, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(92, 5).Value = "SVALÖV"
$ws.Cells.Item(92, 7).Value = 0.8
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(92, 14).Value = 0
$ws.Cells.Item(92, 15).Value = 0
$ws.Cells.Item(92, 16).Value = 0
$ws.Cells.Item(92, 17).Value = 0

# Row 93: A 11633-2024
$ws.Cells.Item(93, 1).Value = "A 11633-2024"
$ws.Cells.Item(93, 2).Value = 45373.4740625
$ws.Cells.Item(93, 3).Value = 46081
$ws.Cells.Item(93, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(93, 5).Value = "SVALÖV"
$ws.Cells.Item(93, 7).Value = 4.9
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(93, 14).Value = 0
$ws.Cells.Item(93, 15).Value = 0
$ws.Cells.Item(93, 16).Value = 0
$ws.Cells.Item(93, 17).Value = 0

# Row 94: A 11922-2023
$ws.Cells.Item(94, 1).Value = "A 11922-2023"
$ws.Cells.Item(94, 2).Value = 44993
$ws.Cells.Item(94, 3).Value = 46081
$ws.Cells.Item(94, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(94, 5).Value = "SVALÖV"
$ws.Cells.Item(94, 7).Value = 0.9
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = 0
$ws.Cells.Item(94, 14).Value = 0
$ws.Cells.Item(94, 15).Value = 0
$ws.Cells.Item(94, 16).Value = 0
$ws.Cells.Item(94, 17).Value = 0

# Row 95: A 14546-2024
$ws.Cells.Item(95, 1).Value = "A 14546-2024"
$ws.Cells.Item(95, 2).Value = 45394
$ws.Cells.Item(95, 3).Value = 46081
$ws.Cells.Item(95, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(95, 5).Value = "SVALÖV"
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 13).Value = 0
$ws.Cells.Item(95, 14).Value = 0
$ws.Cells.Item(95, 15).Value = 0
$ws.Cells.Item(95, 16).Value = 0
$ws.Cells.Item(95, 17).Value = 0

# Row 96: A 32620-2023
$ws.Cells.Item(96, 1).Value = "A 32620-2023"
$ws.Cells.Item(96, 2).Value = 45111
$ws.Cells.Item(96, 3).Value = 46081
$ws.Cells.Item(96, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(96, 5).Value = "SVALÖV"
$ws.Cells.Item(96, 7).Value = 1.4
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = 0
$ws.Cells.Item(96, 14).Value = 0
$ws.Cells.Item(96, 15).Value = 0
$ws.Cells.Item(96, 16).Value = 0
$ws.Cells.Item(96, 17).Value = 0

# Row 97: A 1592-2024
$ws.Cells.Item(97, 1).Value = "A 1592-2024"
$ws.Cells.Item(97, 2).Value = 45306.58594907408
$ws.Cells.Item(97, 3).Value = 46081
$ws.Cells.Item(97, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(97, 5).Value = "SVALÖV"
$ws.Cells.Item(97, 7).Value = 2.2
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(97, 14).Value = 0
$ws.Cells.Item(97, 15).Value = 0
$ws.Cells.Item(97, 16).Value = 0
$ws.Cells.Item(97, 17).Value = 0

# Row 98: A 4149-2022
$ws.Cells.Item(98, 1).Value = "A 4149-2022"
$ws.Cells.Item(98, 2).Value = 44588
$ws.Cells.Item(98, 3).Value = 46081
$ws.Cells.Item(98, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(98, 5).Value = "SVALÖV"
$ws.Cells.Item(98, 7).Value = 1.1
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(98, 14).Value = 0
$ws.Cells.Item(98, 15).Value = 0
$ws.Cells.Item(98, 16).Value = 0
$ws.Cells.Item(98, 17).Value = 0

# Row 99: A 15487-2025
$ws.Cells.Item(99, 1).Value = "A 15487-2025"
$ws.Cells.Item(99, 2).Value = 45747
$ws.Cells.Item(99, 3).Value = 46081
$ws.Cells.Item(99, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(99, 5).Value = "SVALÖV"
$ws.Cells.Item(99, 7).Value = 0.9
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(99, 14).Value = 0
$ws.Cells.Item(99, 15).Value = 0
$ws.Cells.Item(99, 16).Value = 0
$ws.Cells.Item(99, 17).Value = 0

# Row 100: A 11101-2023
$ws.Cells.Item(100, 1).Value = "A 11101-2023"
$ws.Cells.Item(100, 2).Value = 44986
$ws.Cells.Item(100, 3).Value = 46081
$ws.Cells.Item(100, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(100, 5).Value = "SVALÖV"
$ws.Cells.Item(100, 7).Value = 1.5
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = 0
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 15).Value = 0
$ws.Cells.Item(100, 16).Value = 0
$ws.Cells.Item(100, 17).Value = 0

# Row 101: A 5184-2024
$ws.Cells.Item(101, 1).Value = "A 5184-2024"
$ws.Cells.Item(101, 2).Value = 45330
$ws.Cells.Item(101, 3).Value = 46081
$ws.Cells.Item(101, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(101, 5).Value = "SVALÖV"
$ws.Cells.Item(101, 7).Value = 2.2
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = 0
$ws.Cells.Item(101, 14).Value = 0
$ws.Cells.Item(101, 15).Value = 0
$ws.Cells.Item(101, 16).Value = 0
$ws.Cells.Item(101, 17).Value = 0

# Row 102: A 13273-2024
$ws.Cells.Item(102, 1).Value = "A 13273-2024"
$ws.Cells.Item(102, 2).Value = 45386
$ws.Cells.Item(102, 3).Value = 46081
$ws.Cells.Item(102, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(102, 5).Value = "SVALÖV"
$ws.Cells.Item(102, 7).Value = 12.4
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = 0
$ws.Cells.Item(102, 14).Value = 0
$ws.Cells.Item(102, 15).Value = 0
$ws.Cells.Item(102, 16).Value = 0
$ws.Cells.Item(102, 17).Value = 0

# Row 103: A 44196-2024
$ws.Cells.Item(103, 1).Value = "A 44196-2024"
$ws.Cells.Item(103, 2).Value = 45573.318032407406
$ws.Cells.Item(103, 3).Value = 46081
$ws.Cells.Item(103, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(103, 5).Value = "SVALÖV"
$ws.Cells.Item(103, 7).Value = 1.8
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = 0
$ws.Cells.Item(103, 14).Value = 0
$ws.Cells.Item(103, 15).Value = 0
$ws.Cells.Item(103, 16).Value = 0
$ws.Cells.Item(103, 17).Value = 0

# Row 104: A 21507-2025
$ws.Cells.Item(104, 1).Value = "A 21507-2025"
$ws.Cells.Item(104, 2).Value = 45782.61987268519
$ws.Cells.Item(104, 3).Value = 46081
$ws.Cells.Item(104, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(104, 5).Value = "SVALÖV"
$ws.Cells.Item(104, 7).Value = 1.7
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 13).Value = 0
$ws.Cells.Item(104, 14).Value = 0
$ws.Cells.Item(104, 15).Value = 0
$ws.Cells.Item(104, 16).Value = 0
$ws.Cells.Item(104, 17).Value = 0

# Row 105: A 21556-2025
$ws.Cells.Item(105, 1).Value = "A 21556-2025"
$ws.Cells.Item(105, 2).Value = 45782.67826388889
$ws.Cells.Item(105, 3).Value = 46081
$ws.Cells.Item(105, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(105, 5).Value = "SVALÖV"
$ws.Cells.Item(105, 7).Value = 1.3
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 0
$ws.Cells.Item(105, 14).Value = 0
$ws.Cells.Item(105, 15).Value = 0
$ws.Cells.Item(105, 16).Value = 0
$ws.Cells.Item(105, 17).Value = 0

# Row 106: A 21541-2025
$ws.Cells.Item(106, 1).Value = "A 21541-2025"
$ws.Cells.Item(106, 2).Value = 45782.66538194445
$ws.Cells.Item(106, 3).Value = 46081
$ws.Cells.Item(106, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(106, 5).Value = "SVALÖV"
$ws.Cells.Item(106, 7).Value = 3.1
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = 0
$ws.Cells.Item(106, 14).Value = 0
$ws.Cells.Item(106, 15).Value = 0
$ws.Cells.Item(106, 16).Value = 0
$ws.Cells.Item(106, 17).Value = 0

# Row 107: A 21814-2025
$ws.Cells.Item(107, 1).Value = "A 21814-2025"
$ws.Cells.Item(107, 2).Value = 45783.72201388889
$ws.Cells.Item(107, 3).Value = 46081
$ws.Cells.Item(107, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(107, 5).Value = "SVALÖV"
$ws.Cells.Item(107, 7).Value = 2.4
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 0
$ws.Cells.Item(107, 14).Value = 0
$ws.Cells.Item(107, 15).Value = 0
$ws.Cells.Item(107, 16).Value = 0
$ws.Cells.Item(107, 17).Value = 0

# Row 108: A 22203-2025
$ws.Cells.Item(108, 1).Value = "A 22203-2025"
$ws.Cells.Item(108, 2).Value = 45785.653819444444
$ws.Cells.Item(108, 3).Value = 46081
$ws.Cells.Item(108, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(108, 5).Value = "SVALÖV"
$ws.Cells.Item(108, 7).Value = 2.8
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = 0
$ws.Cells.Item(108, 14).Value = 0
$ws.Cells.Item(108, 15).Value = 0
$ws.Cells.Item(108, 16).Value = 0
$ws.Cells.Item(108, 17).Value = 0

# Row 109: A 20047-2024
$ws.Cells.Item(109, 1).Value = "A 20047-2024"
$ws.Cells.Item(109, 2).Value = 45434.373761574076
$ws.Cells.Item(109, 3).Value = 46081
$ws.Cells.Item(109, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(109, 5).Value = "SVALÖV"
$ws.Cells.Item(109, 7).Value = 2.3
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = 0
$ws.Cells.Item(109, 14).Value = 0
$ws.Cells.Item(109, 15).Value = 0
$ws.Cells.Item(109, 16).Value = 0
$ws.Cells.Item(109, 17).Value = 0

# Row 110: A 34808-2024
$ws.Cells.Item(110, 1).Value = "A 34808-2024"
$ws.Cells.Item(110, 2).Value = 45526.69388888889
$ws.Cells.Item(110, 3).Value = 46081
$ws.Cells.Item(110, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(110, 5).Value = "SVALÖV"
$ws.Cells.Item(110, 7).Value = 3.6
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 0
$ws.Cells.Item(110, 14).Value = 0
$ws.Cells.Item(110, 15).Value = 0
$ws.Cells.Item(110, 16).Value = 0
$ws.Cells.Item(110, 17).Value = 0

# Row 111: A 5167-2024
$ws.Cells.Item(111, 1).Value = "A 5167-2024"
$ws.Cells.Item(111, 2).Value = 45330
$ws.Cells.Item(111, 3).Value = 46081
$ws.Cells.Item(111, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(111, 5).Value = "SVALÖV"
$ws.Cells.Item(111, 7).Value = 6.5
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = 0
$ws.Cells.Item(111, 14).Value = 0
$ws.Cells.Item(111, 15).Value = 0
$ws.Cells.Item(111, 16).Value = 0
$ws.Cells.Item(111, 17).Value = 0

# Row 112: A 23747-2025
$ws.Cells.Item(112, 1).Value = "A 23747-2025"
$ws.Cells.Item(112, 2).Value = 45793.44744212963
$ws.Cells.Item(112, 3).Value = 46081
$ws.Cells.Item(112, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(112, 5).Value = "SVALÖV"
$ws.Cells.Item(112, 7).Value = 2.6
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = 0
$ws.Cells.Item(112, 14).Value = 0
$ws.Cells.Item(112, 15).Value = 0
$ws.Cells.Item(112, 16).Value = 0
$ws.Cells.Item(112, 17).Value = 0

# Row 113: A 24513-2025
$ws.Cells.Item(113, 1).Value = "A 24513-2025"
$ws.Cells.Item(113, 2).Value = 45798.40712962963
$ws.Cells.Item(113, 3).Value = 46081
$ws.Cells.Item(113, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(113, 5).Value = "SVALÖV"
$ws.Cells.Item(113, 7).Value = 0.4
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 0
$ws.Cells.Item(113, 14).Value = 0
$ws.Cells.Item(113, 15).Value = 0
$ws.Cells.Item(113, 16).Value = 0
$ws.Cells.Item(113, 17).Value = 0

# Row 114: A 27310-2023
$ws.Cells.Item(114, 1).Value = "A 27310-2023"
$ws.Cells.Item(114, 2).Value = 45096.655497685184
$ws.Cells.Item(114, 3).Value = 46081
$ws.Cells.Item(114, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(114, 5).Value = "SVALÖV"
$ws.Cells.Item(114, 7).Value = 2.5
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 0
$ws.Cells.Item(114, 14).Value = 0
$ws.Cells.Item(114, 15).Value = 0
$ws.Cells.Item(114, 16).Value = 0
$ws.Cells.Item(114, 17).Value = 0
